$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.004000000000001
$ws.Range("B6").Value = 5.708
$ws.Range("B7").Value = 5.396000000000001
$ws.Range("B16").Value = 5.464
$ws.Range("B20").Value = 7.331
$ws.Range("B28").Value = 5.624000000000001
$ws.Range("B29").Value = 5.306
$ws.Range("B32").Value = 6.77
$ws.Range("B40").Value = 9.370000000000001
$ws.Range("B46").Value = 6.382000000000001
$ws.Range("B51").Value = 5.547999999999999
$ws.Range("B52").Value = 5.516
$ws.Range("B57").Value = 5.207
$ws.Range("B59").Value = 4.679
$ws.Range("B62").Value = 5.386
$ws.Range("B66").Value = 5.013
$ws.Range("B73").Value = 7.205000000000001
$ws.Range("B74").Value = 9.164999999999999
$ws.Range("B92").Value = 4.892999999999999
$ws.Range("B100").Value = 5.751
